# Update "想去人数" (interested-count) figures in the 广州-漫展信息 workbook
# to match freshly regenerated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F15").Value = 1049
$ws1.Range("F16").Value = 12179
$ws1.Range("F26").Value = 500
$ws1.Range("F29").Value = 310
$ws1.Range("F31").Value = 274
$ws1.Range("F32").Value = 90
$ws1.Range("F37").Value = 1211

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 329
$ws2.Range("F21").Value = 7

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 1049
$ws4.Range("F14").Value = 12179
$ws4.Range("F24").Value = 500
$ws4.Range("F32").Value = 329
$ws4.Range("F33").Value = 310
$ws4.Range("F36").Value = 274
$ws4.Range("F37").Value = 90
$ws4.Range("F45").Value = 1211
$ws4.Range("F47").Value = 7
